$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 263; existing rows 263-314 shift down to 264-315.
$ws.Rows.Item(263).Insert()

# The new row 263 keeps the same Market/Category/Variety/Quality/Volume/Unit/Origin/Kg
# values as the (now shifted) row 264, but carries a new date and new price figures.
$ws.Cells.Item(263, 1).Value = 9
$ws.Cells.Item(263, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(263, 3).Value = "Metropolitana"
$ws.Cells.Item(263, 4).Value = 45005
$ws.Cells.Item(263, 5).Value = 13
$ws.Cells.Item(263, 6).Value = 100112003
$ws.Cells.Item(263, 7).Value = "Ajo"
$ws.Cells.Item(263, 8).Value = "Chino"
$ws.Cells.Item(263, 9).Value = "Primera"
$ws.Cells.Item(263, 10).Value = 520
$ws.Cells.Item(263, 11).Value = 13000
$ws.Cells.Item(263, 12).Value = 14000
$ws.Cells.Item(263, 13).Value = 13500
$ws.Cells.Item(263, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(263, 15).Value = "China"
$ws.Cells.Item(263, 16).Value = 1350
$ws.Cells.Item(263, 17).Value = 10
$ws.Cells.Item(263, 18).Value = "Hortaliza"
